# "reverting to 3 basic shapes" — undo the addition of the 4th slide
# (the cloud-callout "shape" slide + its notes page), restore the
# triangle's fill from plain yellow to the theme accent1 color, and
# tweak the slide-3 speaker notes wording.

$p = $ppt.ActivePresentation

# 1) Remove the 4th slide (cloud callout) and its notes page along with it.
$p.Slides.Item(4).Delete()

# 2) Triangle on slide 3: solid yellow (FFFF00) -> theme accent1.
$triangle = $p.Slides.Item(3).Shapes.Item(3)
$triangle.Fill.ForeColor.ObjectThemeColor = 5   # msoThemeColorAccent1

# 3) Speaker notes on slide 3: drop the word "yellow".
$notesRange = $p.Slides.Item(3).NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesRange.Text = "Did you predict the next shape would be a triangle?"
